$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the date-format style from A7 to the two new date cells (A8, A9)
$ws.Cells.Item(7,1).Copy()
$ws.Cells.Item(8,1).PasteSpecial(-4122)
$ws.Cells.Item(9,1).PasteSpecial(-4122)

# Row 8: new work entry
$ws.Cells.Item(8,1).Value = 42886
$ws.Cells.Item(8,2).Value = 2.5
$ws.Cells.Item(8,3).Value = "Index - rozvržení, texty, rozvržení bloků apod.."

# Row 9: date only, rest blank
$ws.Cells.Item(9,1).Value = 42887

# Move selection to A9 (also clears the previous scroll/topLeftCell state)
$ws.Range("A9").Select()
